$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BadEvents")

# --- Row 11: replace the old 2-option "bored / rainbow" event with a new
#     5-column "Nothing incredible has happened" bad event (Event + Option 1-4) ---
$ws.Range("A11").Value = "Nothing incredible has happened lately. Your people are bored and complaining, one of them is kicking a tumble weed."
$ws.Range("B11").Value = "Using your brush of magnificance, you paint the sky with a marvelous rainbow. Your people are amused."
$ws.Range("C11").Value = "In an act of generosity, you drip water from your chalice upon the village crops. The villagers rejoice in their harvest."
$ws.Range("D11").Value = "You strike down one of your followers with a brilliant flash of lightning. They certainly aren't bored anymore."
$ws.Range("E11").Value = "In a fit of anger, you rain fire down upon the villagers! Are they not entertained?"
$ws.Rows.Item(11).RowHeight = 30

# --- Row 34: add the new Option 1-4 Good/Bad outcomes for the "nay-sayer"
#     movement event (column A already holds the event description) ---
$ws.Range("B34").Value = "You appear before the nay-sayer in a dream in all your glory… then smack him across the face. The next day he hastily ended the movement."
$ws.Range("C34").Value = "You beset the nay-sayer with nightmares and visions of inconceivable horror. Nothing has changed, he must have forgot his dreams."
$ws.Range("D34").Value = "You give the little heretic a plague of severe itchyness. He seems to have learned his lesson, and now everyday he does extra Holy Jumping Jacks."
$ws.Range("E34").Value = "You plagued the little heretic with relentless sneezing. He seems to have become angry and the movement continues."
$ws.Range("F34").Value = "You struck down the sinful rebel with a blinding flash of lightning.None of your follows skipped their Holy Jumping Jacks that day."
$ws.Range("G34").Value = "In your rage you hurled down a mighty thunderbolt to smite the sinful rebel, but you missed and hit the cheese-maker. The movement doubles in strength."
$ws.Range("H34").Value = "You unleash a vicious plague of mosquitos upon your people until they repent their sinful ways. The nay-sayer quickly loses all his support."
$ws.Range("I34").Value = "You reach to the sun and scoop a mighty fire-ball then drop it upon the rebelious movement. Many people died, but most of them were rebels… you know… probably."
$ws.Rows.Item(34).RowHeight = 60

# --- Update the saved view state: scrolled to A4, active cell now A39 ---
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A39").Select() | Out-Null
